$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("F15").Value = "-"
